# Apply the "I0 and IF added" change: add two new columns I (I0) and J (IF)
# with header labels in row 1 (styled like the existing header row) and
# numeric data in rows 2-73.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header style (bold, bordered, centered) from H1 onto
# the two new header cells, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for rows 2..73 (row r -> index r-2 in the arrays below)
$iVals = @(9,8,9,9,9,9,8,8,9,9,9,9,9,8,9,9,9,9,8,9,8,8,8,9,8,8,9,8,8,9,9,9,9,9,9,9,9,8,9,9,9,9,8,8,8,8,8,9,9,8,8,8,9,9,9,8,8,9,7,9,7,7,9,8,7,9,8,8,5,4,9,9)
$jVals = @(9,8,9,9,9,9,9,8,9,9,9,9,9,8,9,9,9,10,9,9,8,8,8,9,9,8,9,8,9,9,9,9,9,9,9,9,9,9,9,9,9,9,8,9,9,8,9,9,9,8,9,8,9,10,9,9,8,9,8,9,7,8,9,8,8,9,8,9,5,4,9,9)

for ($r = 2; $r -le 73; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iVals[$idx]
    $ws.Cells.Item($r, 10).Value = $jVals[$idx]
}
